# Scheduled runner: refresh market-price-derived profit figures across all Mateus_Profits job sheets.
# Values come from an external price-data pull; no formulas are involved, so every cell is written directly.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 1500
$ws.Range("I18").Value = 1500
$ws.Range("K18").Value = 1500
$ws.Range("M18").Value = -1216
$ws.Range("H131").Value = 112792
$ws.Range("I131").Value = 126578.5
$ws.Range("K131").Value = 379735.5
$ws.Range("M131").Value = -374695.5
$ws.Range("H132").Value = 981.0454999999999
$ws.Range("I132").Value = 879.7105
$ws.Range("J132").Value = 1622.8334
$ws.Range("K132").Value = 2639.1315
$ws.Range("L132").Value = 4868.5002
$ws.Range("M132").Value = -109.1315
$ws.Range("N132").Value = -9928.5002
$ws.Range("H135").Value = 505.65384
$ws.Range("I135").Value = 505.65384
$ws.Range("K135").Value = 4550.88456
$ws.Range("M135").Value = -2015.88456
$ws.Range("H136").Value = 96142.42999999999
$ws.Range("J136").Value = 96142.42999999999
$ws.Range("L136").Value = 96142.42999999999
$ws.Range("N136").Value = -106342.43
$ws.Range("H137").Value = 1596.5454
$ws.Range("J137").Value = 2110
$ws.Range("L137").Value = 6330
$ws.Range("N137").Value = -11430
$ws.Range("H138").Value = 2151.6667
$ws.Range("I138").Value = 1656.4667
$ws.Range("J138").Value = 2342.1282
$ws.Range("K138").Value = 4969.4001
$ws.Range("L138").Value = 7026.3846
$ws.Range("M138").Value = 170.5999000000002
$ws.Range("N138").Value = -17306.3846

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2198.585
$ws.Range("I32").Value = 2251.5881
$ws.Range("K32").Value = 2251.5881
$ws.Range("M32").Value = -1964.5881
$ws.Range("H45").Value = 3778.5
$ws.Range("J45").Value = 6650
$ws.Range("L45").Value = 6650
$ws.Range("N45").Value = -7404
$ws.Range("H61").Value = 7772.1665
$ws.Range("I61").Value = 6660.1924
$ws.Range("K61").Value = 6660.1924
$ws.Range("M61").Value = -6448.1924
$ws.Range("H74").Value = 3752.875
$ws.Range("I74").Value = 2860.68
$ws.Range("J74").Value = 5239.8667
$ws.Range("K74").Value = 2860.68
$ws.Range("L74").Value = 5239.8667
$ws.Range("M74").Value = -1986.68
$ws.Range("N74").Value = -6987.8667
$ws.Range("H77").Value = 3752.875
$ws.Range("I77").Value = 2860.68
$ws.Range("J77").Value = 5239.8667
$ws.Range("K77").Value = 14303.4
$ws.Range("L77").Value = 26199.3335
$ws.Range("M77").Value = -9935.4
$ws.Range("N77").Value = -34935.33349999999
$ws.Range("H97").Value = 533.25
$ws.Range("I97").Value = 631.1111
$ws.Range("J97").Value = 239.66667
$ws.Range("K97").Value = 631.1111
$ws.Range("L97").Value = 239.66667
$ws.Range("M97").Value = -135.1111
$ws.Range("N97").Value = -1231.66667
$ws.Range("H132").Value = 4908.896
$ws.Range("I132").Value = 4045.05
$ws.Range("K132").Value = 12135.15
$ws.Range("M132").Value = -9605.150000000001
$ws.Range("H135").Value = 132000
$ws.Range("J135").Value = 132000
$ws.Range("L135").Value = 132000
$ws.Range("N135").Value = -142140
$ws.Range("H136").Value = 7772.1665
$ws.Range("I136").Value = 6660.1924
$ws.Range("K136").Value = 19980.5772
$ws.Range("M136").Value = -17430.5772

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 5483.9473
$ws.Range("I99").Value = 4154.091
$ws.Range("J99").Value = 7312.5
$ws.Range("K99").Value = 4154.091
$ws.Range("L99").Value = 7312.5
$ws.Range("M99").Value = -2656.091
$ws.Range("N99").Value = -10308.5
$ws.Range("H106").Value = 300000
$ws.Range("I106").Value = 500000
$ws.Range("J106").Value = 100000
$ws.Range("K106").Value = 500000
$ws.Range("L106").Value = 100000
$ws.Range("M106").Value = -498738
$ws.Range("N106").Value = -102524
$ws.Range("H133").Value = 199994.5
$ws.Range("J133").Value = 199994.5
$ws.Range("L133").Value = 199994.5
$ws.Range("N133").Value = -210114.5
$ws.Range("H134").Value = 2672.4902
$ws.Range("I134").Value = 2703.94
$ws.Range("K134").Value = 8111.82
$ws.Range("M134").Value = -5576.82
$ws.Range("H141").Value = 0
$ws.Range("I141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("M141").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3626.7354
$ws.Range("I31").Value = 2600.5217
$ws.Range("K31").Value = 2600.5217
$ws.Range("M31").Value = -2305.5217
$ws.Range("H34").Value = 3626.7354
$ws.Range("I34").Value = 2600.5217
$ws.Range("K34").Value = 2600.5217
$ws.Range("M34").Value = -2398.5217
$ws.Range("H50").Value = 26500
$ws.Range("J50").Value = 26500
$ws.Range("L50").Value = 26500
$ws.Range("N50").Value = -27750
$ws.Range("H58").Value = 7672.5454
$ws.Range("I58").Value = 5386.3335
$ws.Range("J58").Value = 12571.571
$ws.Range("K58").Value = 5386.3335
$ws.Range("L58").Value = 12571.571
$ws.Range("M58").Value = -5183.3335
$ws.Range("N58").Value = -12977.571
$ws.Range("H80").Value = 36957.39
$ws.Range("J80").Value = 36957.39
$ws.Range("L80").Value = 36957.39
$ws.Range("N80").Value = -39203.39
$ws.Range("H83").Value = 36957.39
$ws.Range("J83").Value = 36957.39
$ws.Range("L83").Value = 110872.17
$ws.Range("N83").Value = -122104.17
$ws.Range("H98").Value = 49749.75
$ws.Range("I98").Value = 48999
$ws.Range("K98").Value = 48999
$ws.Range("M98").Value = -46753
$ws.Range("H99").Value = 5221.5454
$ws.Range("I99").Value = 4683.857
$ws.Range("K99").Value = 4683.857
$ws.Range("M99").Value = -3185.857
$ws.Range("H126").Value = 5221.5454
$ws.Range("I126").Value = 4683.857
$ws.Range("K126").Value = 14051.571
$ws.Range("M126").Value = -11581.571
$ws.Range("H132").Value = 2450.081
$ws.Range("I132").Value = 2034
$ws.Range("K132").Value = 6102
$ws.Range("M132").Value = -3572
$ws.Range("H134").Value = 5388.032
$ws.Range("I134").Value = 4400.6
$ws.Range("K134").Value = 13201.8
$ws.Range("M134").Value = -10666.8
$ws.Range("H136").Value = 7672.5454
$ws.Range("I136").Value = 5386.3335
$ws.Range("J136").Value = 12571.571
$ws.Range("K136").Value = 16159.0005
$ws.Range("L136").Value = 37714.713
$ws.Range("M136").Value = -13609.0005
$ws.Range("N136").Value = -42814.713

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 778.5
$ws.Range("I113").Value = 685.5
$ws.Range("J113").Value = 815.7
$ws.Range("K113").Value = 2056.5
$ws.Range("L113").Value = 2447.1
$ws.Range("M113").Value = 113.5
$ws.Range("N113").Value = -6787.1

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 9159.182000000001
$ws.Range("I2").Value = 75.2
$ws.Range("K2").Value = 75.2
$ws.Range("M2").Value = 37.8
$ws.Range("H36").Value = 5017
$ws.Range("I36").Value = 5017
$ws.Range("K36").Value = 5017
$ws.Range("M36").Value = -4532
$ws.Range("H55").Value = 31662.666
$ws.Range("I55").Value = 15000
$ws.Range("K55").Value = 15000
$ws.Range("M55").Value = -14673
$ws.Range("H122").Value = 3200.3333
$ws.Range("I122").Value = 3257.0588
$ws.Range("J122").Value = 3062.5715
$ws.Range("K122").Value = 9771.1764
$ws.Range("L122").Value = 9187.7145
$ws.Range("M122").Value = -7321.1764
$ws.Range("N122").Value = -14087.7145
$ws.Range("H126").Value = 4685.852
$ws.Range("J126").Value = 5074
$ws.Range("L126").Value = 15222
$ws.Range("N126").Value = -20162
$ws.Range("H132").Value = 1635.3784
$ws.Range("I132").Value = 1635.3784
$ws.Range("K132").Value = 4906.135200000001
$ws.Range("M132").Value = -2376.135200000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 10877.069
$ws.Range("I132").Value = 12620.4
$ws.Range("J132").Value = 7003
$ws.Range("K132").Value = 37861.2
$ws.Range("L132").Value = 21009
$ws.Range("M132").Value = -35331.2
$ws.Range("N132").Value = -26069
$ws.Range("H136").Value = 4510.32
$ws.Range("I136").Value = 4088.8096
$ws.Range("K136").Value = 12266.4288
$ws.Range("M136").Value = -9716.4288
$ws.Range("H137").Value = 75100
$ws.Range("J137").Value = 75100
$ws.Range("L137").Value = 75100
$ws.Range("N137").Value = -85300
$ws.Range("H139").Value = 82699.5
$ws.Range("J139").Value = 82699.5
$ws.Range("L139").Value = 82699.5
$ws.Range("N139").Value = -92979.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 5331.067
$ws.Range("I62").Value = 4904.273
$ws.Range("K62").Value = 4904.273
$ws.Range("M62").Value = -4280.273
$ws.Range("H65").Value = 5331.067
$ws.Range("I65").Value = 4904.273
$ws.Range("K65").Value = 24521.365
$ws.Range("M65").Value = -21401.365
$ws.Range("H96").Value = 2412.182
$ws.Range("I96").Value = 1407
$ws.Range("J96").Value = 3249.8333
$ws.Range("K96").Value = 1407
$ws.Range("L96").Value = 3249.8333
$ws.Range("M96").Value = -34
$ws.Range("N96").Value = -5995.8333
$ws.Range("H122").Value = 4347.9585
$ws.Range("I122").Value = 3929.3684
$ws.Range("K122").Value = 11788.1052
$ws.Range("M122").Value = -9338.1052
$ws.Range("H136").Value = 4217.3213
$ws.Range("I136").Value = 2448.7727
$ws.Range("J136").Value = 10702
$ws.Range("K136").Value = 7346.3181
$ws.Range("L136").Value = 32106
$ws.Range("M136").Value = -4796.3181
$ws.Range("N136").Value = -37206

Write-Output "Applied 238 cell updates across $($wb.Worksheets.Count) sheets"
